$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 400
$ws.Range("B3").Value = 300
$ws.Range("B5").Value = 57
$ws.Range("B6").Value = 160
$ws.Range("B7").Value = 450
$ws.Range("B8").Value = 535
